$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted before the current row 114, pushing the
# existing rows 114-136 down to 115-137 (their data stays unchanged).
$ws.Rows.Item(114).Insert()

# Populate the newly inserted row 114 with the new record's data.
$ws.Cells.Item(114, 1).Value = 5
$ws.Cells.Item(114, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(114, 3).Value = "Maule"
$ws.Cells.Item(114, 4).Value = 44504
$ws.Cells.Item(114, 5).Value = 7
$ws.Cells.Item(114, 6).Value = 100112017
$ws.Cells.Item(114, 7).Value = "Apio"
$ws.Cells.Item(114, 8).Value = "Americana (o)"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 500
$ws.Cells.Item(114, 11).Value = 7500
$ws.Cells.Item(114, 12).Value = 7500
$ws.Cells.Item(114, 13).Value = 7500
$ws.Cells.Item(114, 14).Value = "$/docena de matas"
$ws.Cells.Item(114, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(114, 16).Value = 1250
$ws.Cells.Item(114, 17).Value = 6
$ws.Cells.Item(114, 18).Value = "Hortaliza"
